$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the run ", LIMIT, and ORDER BY" that used to follow
# "WHERE" in the assignment title ("... Select data with WHERE clause.").
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(", LIMIT, and ORDER BY")
if ($found) {
    $startPos = $rng.Start
    $rng.Text = ""

    # Deleting that run merges the formatting-identical runs that follow it
    # in the same paragraph (the space / "clause" / "." runs) into a single
    # run. Nudge the Bold flag off and back on (a net no-op - the paragraph
    # is already bold) on each of those pieces so the engine re-splits them
    # back into separate runs, matching the original (untouched) structure.
    $space = $d.Range($startPos, $startPos + 1)
    if ($space.Text -eq " ") {
        $space.Font.Bold = 0
        $space.Font.Bold = 1
    }
    $clause = $d.Range($startPos + 1, $startPos + 7)
    if ($clause.Text -eq "clause") {
        $clause.Font.Bold = 0
        $clause.Font.Bold = 1
    }
}

# ---------------------------------------------------------------------------
# Change 2: "Display all employee information of the employee ID is 15."
# -> "Display employee information of the employee ID is 15."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Display all employee information of the employee", $true, $false, $false, $false, $false, $true, 1, $false, "Display employee information of the employee", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: merge the two runs "Display all employee " + "whose salary is
# more than 4000." into a single run/sentence (text unchanged, only the run
# split goes away).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Display all employee whose salary is more than 4000.", $true, $false, $false, $false, $false, $true, 1, $false, "Display all employee whose salary is more than 4000.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: merge the three runs "Display all employee whose salary is more
# than 4000" + " and less than 5000" + "." into a single run/sentence (text
# unchanged, only the run split goes away).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Display all employee whose salary is more than 4000 and less than 5000.", $true, $false, $false, $false, $false, $true, 1, $false, "Display all employee whose salary is more than 4000 and less than 5000.", 2) | Out-Null
